$wb = $excel.ActiveWorkbook

# --- Sheet "SoCDTtiNTY-psgr" ---
$ws1 = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$ws1.Range("D7").Value = 0.068
$ws1.Range("D8").Select()

# --- Sheet "SoCDTtiNTY-frgt" ---
$ws2 = $wb.Worksheets.Item("SoCDTtiNTY-frgt")
$ws2.Range("B4:H4").Value = 0.028
$ws2.Range("B4:H4").Select()

# Restore the originally active sheet ("About") so the workbook re-opens
# on the same tab as before the edit.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
